# AutoCommit_6 июня 2024 г. 14:01:13_SibNout2023
# Adds a new "МЫШАКОВ!!!" summary row directly below the last student row
# (row 31), giving it five "5" scores formatted like the other data rows,
# and nudges the view down to the new bottom of the sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New row 31: label + five "5" scores ---
$ws.Range("B31").Value = "МЫШАКОВ!!!"
$ws.Range("C31:G31").Value = 5

# Style B31 like the other bold/bordered "name" column header cells
# (bold font, centered, wrapped, thick border) and C31:G31 like the
# thick-bordered score cells used throughout the sheet.
$ws.Range("B31").Font.Bold = $true
$ws.Range("B31").HorizontalAlignment = -4108  # xlCenter
$ws.Range("B31").VerticalAlignment = -4108    # xlCenter
$ws.Range("B31").WrapText = $true

$ws.Range("C31:G31").HorizontalAlignment = -4108  # xlCenter
$ws.Range("C31:G31").VerticalAlignment = -4108    # xlCenter
$ws.Range("C31:G31").WrapText = $true
$ws.Range("C31:G31").Font.Bold = $false
$ws.Range("C31:G31").Font.Color = 0
$ws.Range("C31:G31").Interior.Pattern = -4142  # xlNone

$thick = -4118  # xlThick
1, 2, 3, 4 | ForEach-Object {
    $ws.Range("C31:G31").Borders.Item($_).LineStyle = 1
    $ws.Range("C31:G31").Borders.Item($_).Weight = $thick
    $ws.Range("C31:G31").Borders.Item($_).Color = 0
}

# --- View: scroll/select the new bottom of the sheet ---
$ws.Range("H31").Select()
